$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark status column "C" as SI/PROBLEMAS for the first few tasks
$ws.Range("C2").Value = "SI"
$ws.Range("C4").Value = "SI"
$ws.Range("C5").Value = "SI"
$ws.Range("C9").Value = "PROBLEMAS"
$ws.Range("C10").Value = "PROBLEMAS"
$ws.Range("C11").Value = "SI"

# Move the active selection to B17 (work in progress)
$ws.Range("B17").Select()
